# Commit: "plotted phylum level t8 lca"
#
# On the "plotting in python" sheet (the data feeding the phylum-level T8
# LCA plot), a third data column (C, "Time 0 small surf") is filled in with
# 0 for every phylum row that did not already have a value there, so the
# whole A1:C38 block is now fully populated and ready to be charted/plotted.
# Two sheets' active-cell selections also moved, reflecting where the user
# was working afterwards.

$wb = $excel.ActiveWorkbook

# --- "plotting in python": fill in column C (value 0) for every row that ---
# --- was still missing it, so the range becomes a complete 3-column table ---
$ws3 = $wb.Worksheets.Item("plotting in python")

$rowsNeedingC = @(4, 6, 7, 8, 9, 10, 11, 12, 13, 14, 15, 16, 17, 18, 19, 20, 21, 22, 24, 25, 26, 27, 28, 29, 30, 31, 32, 33, 34, 35, 36, 37, 38)
foreach ($r in $rowsNeedingC) {
    $ws3.Range("C$r").Value = 0
}

# Scroll the view down (top row 19) and move the selection to D38.
$ws3.Activate()
$excel.ActiveWindow.ScrollRow = 19
$ws3.Range("D38").Select()

# --- "sorting 675_50": active cell moved from D1 to C2 ---
$ws2 = $wb.Worksheets.Item("sorting 675_50")
$ws2.Activate()
$ws2.Range("C2").Select()

# "plotting in python" remains the active/selected tab, with its new
# selection (D38) in place.
$ws3.Activate()
$ws3.Range("D38").Select()
